# "Vorbereitung für pdf als medium"
# - Rename the sheet "Tabelle1" to "ProbeDatenPDF" (it will hold the data
#   that feeds the PDF export).
# - Make that sheet the active/selected tab (previously "Kontrolle" was
#   the active tab).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Name = "ProbeDatenPDF"

# Switch the active tab from "Kontrolle" to the renamed sheet - this also
# flips sheetView/tabSelected between the two worksheets and updates
# bookViews/workbookView@activeTab accordingly.
$ws.Activate()
